# Apply the "7/3-2018 Wave collision" diary update to the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new laboratory diary entries (rows 19-21) -------------------
# Row 19: 2/3-2018 - Projectile implemented - 0h 40m
$ws.Cells.Item(19, 1).Value = "2/3-2018"
$ws.Cells.Item(19, 2).Value = "Projectile implemented"
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 40

# Row 20: 2/3-2018 - Collision implemented - 1h 0m
$ws.Cells.Item(20, 1).Value = "2/3-2018"
$ws.Cells.Item(20, 2).Value = "Collision implemented"
$ws.Cells.Item(20, 3).Value = 1
$ws.Cells.Item(20, 4).Value = 0

# Row 21: 7/3-2018 - Wave collision - 4h 30m
$ws.Cells.Item(21, 1).Value = "7/3-2018"
$ws.Cells.Item(21, 2).Value = "Wave collision"
$ws.Cells.Item(21, 3).Value = 4
$ws.Cells.Item(21, 4).Value = 30

# --- Update the selection to reflect where the user left off -------------
$ws.Range("F21").Select()

# --- Recalculate totals (formulas already reference C2:C21 / D2:D21) -----
$excel.Calculate()

# --- Reflect the new window size Excel recorded for the workbook ---------
# (Best-effort; the hosted window size is an OS-level attribute mirrored
# from the active window when Excel saves the file.)
$win = $excel.ActiveWindow
$win.Width = 17256
$win.Height = 5640
